# Add a new "ad_core_geog" baseline table to the BiBBS_Geographic sheet,
# mirroring the existing "ch_core_geog" table, inserted directly above the
# existing "house_nbhd" rows inside Table5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BiBBS_Geographic")

# Make room for 13 new data rows just above the current row 15
# (this shifts the existing house_nbhd rows from 15-20 down to 28-33).
$ws.Rows("15:27").Insert()

# New rows to insert (project, table, variable, full_name, label, value_type,
# description, categories, categories_label).
$newRows = @(
    @("BiBBS_Geographic", "ad_core_geog", "date_address_data", "BiBBS_Geographic.ad_core_geog.date_address_data", "Date of latest available address data", "date", "", "", ""),
    @("BiBBS_Geographic", "ad_core_geog", "age_m", "BiBBS_Geographic.ad_core_geog.age_m", "Participants actual age (months)", "decimal", "", "", ""),
    @("BiBBS_Geographic", "ad_core_geog", "age_closest_data_point", "BiBBS_Geographic.ad_core_geog.age_closest_data_point", "Age at data point closest to participants actual age", "decimal", "", "", ""),
    @("BiBBS_Geographic", "ad_core_geog", "temporal_accuracy_m", "BiBBS_Geographic.ad_core_geog.temporal_accuracy_m", "Difference in months between participants age at closest data point and actual a", "decimal", "", "", ""),
    @("BiBBS_Geographic", "ad_core_geog", "LSOA11CD", "BiBBS_Geographic.ad_core_geog.LSOA11CD", "LSOA 2011 code", "text", "", "", ""),
    @("BiBBS_Geographic", "ad_core_geog", "WD21CD", "BiBBS_Geographic.ad_core_geog.WD21CD", "Ward 2021 code", "text", "", "", ""),
    @("BiBBS_Geographic", "ad_core_geog", "is_in_bfd_la", "BiBBS_Geographic.ad_core_geog.is_in_bfd_la", "Is participant in Bradford LA?", "decimal", "", "0|1", "No|Yes"),
    @("BiBBS_Geographic", "ad_core_geog", "is_in_bibbs_area", "BiBBS_Geographic.ad_core_geog.is_in_bibbs_area", "Is participant in BiBBS area?", "decimal", "", "0|1", "No|Yes"),
    @("BiBBS_Geographic", "ad_core_geog", "data_source", "BiBBS_Geographic.ad_core_geog.data_source", "Source of data: registration (1) or tracing (2)", "decimal", "", "1|2", "Registration|Tracing"),
    @("BiBBS_Geographic", "ad_core_geog", "study", "BiBBS_Geographic.ad_core_geog.study", "Is participant in BiB (1) or BiBBS (2)", "decimal", "", "1|2", "BiB|BiBBS"),
    @("BiBBS_Geographic", "ad_core_geog", "not_in_eng_wales", "BiBBS_Geographic.ad_core_geog.not_in_eng_wales", "Indicates if address is not in England or Wales", "decimal", "", "0|1", "No|Yes"),
    @("BiBBS_Geographic", "ad_core_geog", "missing_address_data", "BiBBS_Geographic.ad_core_geog.missing_address_data", "Indicates if record has missing address data", "decimal", "", "0|1", "No|Yes"),
    @("BiBBS_Geographic", "ad_core_geog", "poor_qual_data", "BiBBS_Geographic.ad_core_geog.poor_qual_data", "Indicates if record has poor quality data", "decimal", "", "0|1", "No|Yes")
)

$startRow = 15
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $val = $rowVals[$c - 1]
        if ($val -ne "") {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# Grow the Table5 listobject so the new rows (and the pre-existing
# house_nbhd rows, now shifted to 28-33) are part of the table again.
$lo = $ws.ListObjects.Item("Table5")
$lo.Resize($ws.Range("A1:I33"))
